# Update experiment with latest results
#
# 1) The "last modified" footer/header date stamp (1/25/2018 -> 2/22/2018)
#    is refreshed everywhere it is cached: the slide master, every slide
#    layout, and the notes master.
# 2) Slide 15's second results table gets a refreshed F1 score for the
#    "Untrained embeddings" row (0.73 -> 0.74).

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes, $newDate) {
    for ($shpIdx = 1; $shpIdx -le $shapes.Count; $shpIdx++) {
        $shp = $shapes.Item($shpIdx)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$newDate = "2/22/2018"

# Slide master.
Update-DatePlaceholder $p.SlideMaster.Shapes $newDate

# Every slide layout hanging off the master.
$layouts = $p.SlideMaster.CustomLayouts
for ($layIdx = 1; $layIdx -le $layouts.Count; $layIdx++) {
    Update-DatePlaceholder $layouts.Item($layIdx).Shapes $newDate
}

# Notes master.
Update-DatePlaceholder $p.NotesMaster.Shapes $newDate

# Slide 15 ("Clone detection experiment results"): refresh the F1 score
# for "Untrained embeddings" from 0.73 to 0.74 in the second table.
$slide = $p.Slides.Item(15)
for ($shIdx = 1; $shIdx -le $slide.Shapes.Count; $shIdx++) {
    $shp = $slide.Shapes.Item($shIdx)
    if ($shp.HasTable) {
        $tbl = $shp.Table
        for ($rowIdx = 1; $rowIdx -le $tbl.Rows.Count; $rowIdx++) {
            $rowLabel = $tbl.Cell($rowIdx, 1).Shape.TextFrame.TextRange.Text
            if ($rowLabel -eq "Untrained embeddings") {
                $f1Cell = $tbl.Cell($rowIdx, 2).Shape.TextFrame.TextRange
                if ($f1Cell.Text -eq "0.73") {
                    $f1Cell.Text = "0.74"
                }
            }
        }
    }
}
